# Update the "need_to_buy" data table with refreshed values (rows 2-15,
# columns A-F). This mirrors a data refresh pulled in from R: dates shift
# forward by one day and the numeric columns get new computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=46042; B=11938.7821073812; C=11733.4550187254; D=17499.86; E=7290.84378325755; F=63.5182834159565}
    @{Row=3;  A=46043; B=12770.1728165669; C=12532.7067109633; D=12075.86; E=8661.76498904513; F=379.94215416702}
    @{Row=4;  A=46044; B=12750.8604180183; C=12533.2166091539; D=12075.86; E=8650.48022336465; F=379.493201354939}
    @{Row=5;  A=46045; B=12620.333747298;  C=11635.2863534695; D=12075.86; E=8538.50748250125; F=337.413909832114}
    @{Row=6;  A=46046; B=4857.99678718736; C=7740.69835097971;  D=12075.86; E=7969.44059015702; F=151.42828921403}
    @{Row=7;  A=46047; B=5000.01932310789; C=7937.34599887101;  D=12075.86; E=8256.71023742689; F=171.591509845746}
    @{Row=8;  A=46048; B=12315.9682835607; C=11865.9738711145; D=12075.86; E=8314.29255877511; F=337.683601245402}
    @{Row=9;  A=46049; B=12315.9682835607; C=12203.7455132662; D=12075.86; E=8314.29255877511; F=351.757419668389}
    @{Row=10; A=46050; B=12315.9682835607; C=12233.2157625311; D=12075.86; E=8314.29255877511; F=352.985346721093}
    @{Row=11; A=46051; B=12315.9682835607; C=11885.4017468464; D=12075.86; E=8314.29255877511; F=338.493096067561}
    @{Row=12; A=46052; B=12315.9682835607; C=11438.0382018787; D=12075.86; E=8314.29255877511; F=319.852948360575}
    @{Row=13; A=46053; B=4867.38022112383; C=8291.09524694406;  D=12075.86; E=7930.23077371803; F=172.72775086092}
    @{Row=14; A=46054; B=5046.61366744637; C=8141.40411370133;  D=9743.86;  E=7890.29102711637; F=261.993130867404}
    @{Row=15; A=46055; B=11445.6138712783; C=11515.348762324;   D=9743.86;  E=7806.00389249187; F=399.062193950659}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
}
